$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1653.75
$ws.Range("I70").Value = 1906.1111
$ws.Range("J70").Value = 896.6667
$ws.Range("K70").Value = 5718.3333
$ws.Range("L70").Value = 2690.0001
$ws.Range("M70").Value = -5448.3333
$ws.Range("N70").Value = -3230.0001
$ws.Range("H73").Value = 1653.75
$ws.Range("I73").Value = 1906.1111
$ws.Range("J73").Value = 896.6667
$ws.Range("K73").Value = 5718.3333
$ws.Range("L73").Value = 2690.0001
$ws.Range("M73").Value = -4782.3333
$ws.Range("N73").Value = -4562.0001
$ws.Range("H127").Value = 2305.3333
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 2305.3333
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 6915.999899999999
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -16835.9999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 20945.2
$ws.Range("I33").Value = 6363
$ws.Range("J33").Value = 30666.666
$ws.Range("K33").Value = 6363
$ws.Range("L33").Value = 30666.666
$ws.Range("M33").Value = -6034
$ws.Range("N33").Value = -31324.666
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H45").Value = 112648.89
$ws.Range("I45").Value = 200966
$ws.Range("J45").Value = 2252.5
$ws.Range("K45").Value = 200966
$ws.Range("L45").Value = 2252.5
$ws.Range("M45").Value = -200589
$ws.Range("N45").Value = -3006.5
$ws.Range("H64").Value = 17586.334
$ws.Range("J64").Value = 17586.334
$ws.Range("L64").Value = 17586.334
$ws.Range("N64").Value = -18082.334
$ws.Range("H67").Value = 17586.334
$ws.Range("J67").Value = 17586.334
$ws.Range("L67").Value = 17586.334
$ws.Range("N67").Value = -19302.334
$ws.Range("H74").Value = 1594.7646
$ws.Range("I74").Value = 1962.4
$ws.Range("J74").Value = 1441.5834
$ws.Range("K74").Value = 1962.4
$ws.Range("L74").Value = 1441.5834
$ws.Range("M74").Value = -1088.4
$ws.Range("N74").Value = -3189.5834
$ws.Range("H77").Value = 1594.7646
$ws.Range("I77").Value = 1962.4
$ws.Range("J77").Value = 1441.5834
$ws.Range("K77").Value = 9812
$ws.Range("L77").Value = 7207.916999999999
$ws.Range("M77").Value = -5444
$ws.Range("N77").Value = -15943.917
$ws.Range("H122").Value = 1694.6666
$ws.Range("I122").Value = 1540.25
$ws.Range("K122").Value = 4620.75
$ws.Range("M122").Value = -2170.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 50000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 50000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 50000
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -50832
$ws.Range("H62").Value = 19176.8
$ws.Range("J62").Value = 19176.8
$ws.Range("L62").Value = 19176.8
$ws.Range("N62").Value = -20548.8
$ws.Range("H65").Value = 19176.8
$ws.Range("J65").Value = 19176.8
$ws.Range("L65").Value = 57530.39999999999
$ws.Range("N65").Value = -64394.39999999999
$ws.Range("H99").Value = 1916.8572
$ws.Range("I99").Value = 969.4
$ws.Range("J99").Value = 2443.2222
$ws.Range("K99").Value = 969.4
$ws.Range("L99").Value = 2443.2222
$ws.Range("M99").Value = 528.6
$ws.Range("N99").Value = -5439.2222
$ws.Range("H107").Value = 41667450
$ws.Range("I107").Value = 45455336
$ws.Range("K107").Value = 45455336
$ws.Range("M107").Value = -45453416
$ws.Range("H134").Value = 2528478.2
$ws.Range("J134").Value = 13902888
$ws.Range("L134").Value = 41708664
$ws.Range("N134").Value = -41713734
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 24233.334
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 28980
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 28980
$ws.Range("M17").Value = -326
$ws.Range("N17").Value = -29328
$ws.Range("H29").Value = 24000
$ws.Range("J29").Value = 24000
$ws.Range("L29").Value = 24000
$ws.Range("N29").Value = -24586
$ws.Range("H99").Value = 142860240
$ws.Range("I99").Value = 200003000
$ws.Range("J99").Value = 3350
$ws.Range("K99").Value = 200003000
$ws.Range("L99").Value = 3350
$ws.Range("M99").Value = -200001502
$ws.Range("N99").Value = -6346
$ws.Range("H100").Value = 30780
$ws.Range("J100").Value = 30780
$ws.Range("L100").Value = 30780
$ws.Range("N100").Value = -32944
$ws.Range("H107").Value = 586.13635
$ws.Range("I107").Value = 453.82352
$ws.Range("K107").Value = 453.82352
$ws.Range("M107").Value = 1466.17648
$ws.Range("H126").Value = 142860240
$ws.Range("I126").Value = 200003000
$ws.Range("J126").Value = 3350
$ws.Range("K126").Value = 600009000
$ws.Range("L126").Value = 10050
$ws.Range("M126").Value = -600006530
$ws.Range("N126").Value = -14990
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 28591.514
$ws.Range("J12").Value = 43506.957
$ws.Range("L12").Value = 130520.871
$ws.Range("N12").Value = -130866.871
$ws.Range("H131").Value = 813.03
$ws.Range("J131").Value = 820.23956
$ws.Range("L131").Value = 2460.71868
$ws.Range("N131").Value = -12540.71868
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 30657
$ws.Range("J101").Value = 30657
$ws.Range("L101").Value = 30657
$ws.Range("N101").Value = -37147
$ws.Range("H102").Value = 2631.1667
$ws.Range("I102").Value = 1193.3334
$ws.Range("J102").Value = 4069
$ws.Range("K102").Value = 1193.3334
$ws.Range("L102").Value = 4069
$ws.Range("M102").Value = 428.6666
$ws.Range("N102").Value = -7313
$ws.Range("H122").Value = 41676484
$ws.Range("I122").Value = 71444104
$ws.Range("J122").Value = 1816
$ws.Range("K122").Value = 214332312
$ws.Range("L122").Value = 5448
$ws.Range("M122").Value = -214329862
$ws.Range("N122").Value = -10348
$ws.Range("H126").Value = 2395.2
$ws.Range("I126").Value = 1988
$ws.Range("J126").Value = 2666.6667
$ws.Range("K126").Value = 5964
$ws.Range("L126").Value = 8000.000100000001
$ws.Range("M126").Value = -3494
$ws.Range("N126").Value = -12940.0001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1779.0714
$ws.Range("I7").Value = 1242.3334
$ws.Range("J7").Value = 4999.5
$ws.Range("K7").Value = 1242.3334
$ws.Range("L7").Value = 4999.5
$ws.Range("M7").Value = -1130.3334
$ws.Range("N7").Value = -5223.5
$ws.Range("H17").Value = 3617.2
$ws.Range("I17").Value = 1396.5
$ws.Range("J17").Value = 12500
$ws.Range("K17").Value = 1396.5
$ws.Range("L17").Value = 12500
$ws.Range("M17").Value = -1226.5
$ws.Range("N17").Value = -12840
$ws.Range("H40").Value = 41670136
$ws.Range("I40").Value = 3002
$ws.Range("J40").Value = 62503700
$ws.Range("K40").Value = 3002
$ws.Range("L40").Value = 62503700
$ws.Range("M40").Value = -2866
$ws.Range("N40").Value = -62503972
$ws.Range("H55").Value = 5749
$ws.Range("J55").Value = 271.66666
$ws.Range("L55").Value = 271.66666
$ws.Range("N55").Value = -617.66666
$ws.Range("H106").Value = 333355140
$ws.Range("J106").Value = 333355140
$ws.Range("L106").Value = 333355140
$ws.Range("N106").Value = -333357664
$ws.Range("H126").Value = 1779.0714
$ws.Range("I126").Value = 1242.3334
$ws.Range("J126").Value = 4999.5
$ws.Range("K126").Value = 3727.0002
$ws.Range("L126").Value = 14998.5
$ws.Range("M126").Value = -1257.0002
$ws.Range("N126").Value = -19938.5
$ws.Range("H132").Value = 40825756
$ws.Range("I132").Value = 76191960
$ws.Range("J132").Value = 18592.385
$ws.Range("K132").Value = 228575880
$ws.Range("L132").Value = 55777.155
$ws.Range("M132").Value = -228573350
$ws.Range("N132").Value = -60837.155
$ws.Range("H136").Value = 119050000
$ws.Range("I136").Value = 76193144
$ws.Range("K136").Value = 228579432
$ws.Range("M136").Value = -228576882
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 100
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 100
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 15
$ws.Range("N6").ClearContents()
$ws.Range("H16").Value = 59800
$ws.Range("J16").Value = 59800
$ws.Range("L16").Value = 59800
$ws.Range("N16").Value = -60384
$ws.Range("H100").Value = 476.55
$ws.Range("I100").Value = 389.4375
$ws.Range("K100").Value = 778.875
$ws.Range("M100").Value = -237.875
$ws.Range("H105").Value = 54480
$ws.Range("J105").Value = 54480
$ws.Range("L105").Value = 54480
$ws.Range("N105").Value = -61468
$ws.Range("H107").Value = 513.3333
$ws.Range("I107").Value = 250
$ws.Range("K107").Value = 750
$ws.Range("M107").Value = 1170
$ws.Range("H122").Value = 48535.727
$ws.Range("I122").Value = 100797.6
$ws.Range("J122").Value = 4984.1665
$ws.Range("K122").Value = 302392.8
$ws.Range("L122").Value = 14952.4995
$ws.Range("M122").Value = -299942.8
$ws.Range("N122").Value = -19852.4995
$ws.Range("H126").Value = 1766
$ws.Range("I126").Value = 828.7143
$ws.Range("J126").Value = 2586.125
$ws.Range("K126").Value = 2486.1429
$ws.Range("L126").Value = 7758.375
$ws.Range("M126").Value = -16.14289999999983
$ws.Range("N126").Value = -12698.375
$ws.Range("H132").Value = 44096.816
$ws.Range("I132").Value = 96137.37
$ws.Range("J132").Value = 8318.9375
$ws.Range("K132").Value = 288412.11
$ws.Range("L132").Value = 24956.8125
$ws.Range("M132").Value = -285882.11
$ws.Range("N132").Value = -30016.8125
$ws.Range("H136").Value = 2117.0715
$ws.Range("I136").Value = 1437.7646
$ws.Range("K136").Value = 4313.293799999999
$ws.Range("M136").Value = -1763.293799999999
